# Updates cryptos list values (prices/volumes/coin reorderings) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.172.04"
$ws.Range("E2").Value = "  -0.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.955.98"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.91"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.87"
$ws.Range("E6").Value = "  -1.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.538"
$ws.Range("E7").Value = "  -0.64%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("E10").Value = "  -1.27%  "

$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0841"
$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.421.71"
$ws.Range("E13").Value = "  +0.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.06"
$ws.Range("E14").Value = "  -2.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.43"
$ws.Range("E15").Value = "  +0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.944.98"
$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.985"
$ws.Range("E17").Value = "  +4.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.131.30"
$ws.Range("E18").Value = "  -0.77%  "

$ws.Range("E19").Value = "  -6.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.11"
$ws.Range("E20").Value = "  -3.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.54"
$ws.Range("E21").Value = "  -3.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0953"
$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.80"
$ws.Range("E24").Value = "  -0.26%  "

$ws.Range("E25").Value = "  +1.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.41"
$ws.Range("E26").Value = "  +14.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.63"
$ws.Range("E27").Value = "  +6.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.169"
$ws.Range("E28").Value = "  +0.28%  "

$ws.Range("E29").Value = "  +10.15%  "

$ws.Range("E30").Value = "  -0.95%  "

$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.72"
$ws.Range("E32").Value = "  -0.64%  "

$ws.Range("E33").Value = "  -0.38%  "

$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0455"
$ws.Range("E34").Value = "  +5.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.91"
$ws.Range("E35").Value = "  -0.94%  "

$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.44"
$ws.Range("E36").Value = "  -2.85%  "

$ws.Range("B37").Value = "Toncoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.05"
$ws.Range("E37").Value = "  -2.08%  "

$ws.Range("E38").Value = "  -0.09%  "

$ws.Range("E39").Value = "  -1.97%  "

$ws.Range("E40").Value = "  -1.29%  "

$ws.Range("E41").Value = "  -1.42%  "

$ws.Range("E42").Value = "  +0.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.78"
$ws.Range("E43").Value = "  -2.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.72"
$ws.Range("E44").Value = "  -2.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.13"
$ws.Range("E45").Value = "  -3.67%  "

$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("E48").Value = "  +2.27%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.007.12"
$ws.Range("E49").Value = "  -0.90%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.23"
$ws.Range("E50").Value = "  +1.49%  "

$ws.Range("E51").Value = "  +3.98%  "
